$wb = $excel.ActiveWorkbook

# The "Spillover estimate" feature needs a new "administrative sanction
# authority" field (mirroring the existing "technical sanction authority"
# field on the technicalSanctionDetails sheet), so add a new column D to
# the adminSanctionDetails sheet with a header + sample value.
$wsAdmin = $wb.Worksheets.Item("adminSanctionDetails")
$wsTech  = $wb.Worksheets.Item("technicalSanctionDetails")

$wsAdmin.Range("D1").Value = "administrativeSanctionAuthority"
$wsAdmin.Range("D2").Value = "commisioner"

# Match the look of the other header cells (bold, green, Arial 9) by
# copying the formatting from the equivalent header cell on the
# technicalSanctionDetails sheet.
$wsTech.Range("B1").Copy()
$wsAdmin.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# adminSanctionDetails becomes the active/selected sheet (previously it
# was technicalSanctionDetails), with D2 as the active cell.
$wsAdmin.Activate()
$wsAdmin.Range("D2").Select()
